$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Số lượng" (H1) renamed to "kích thước" ---
$ws.Range("H1").Value = "kích thước"

# --- Row 2 (STT 1, BTC/USDT trade): note text flipped from "mua..." to "bán..." ---
$ws.Range("J2").Value = "bán khi macd âm, ema âm"

# --- Highlight rows 3, 4 and 6 with red font color (closed/reviewed trades) ---
foreach ($r in 3,4,6) {
    foreach ($col in "A","B","C","D","E","F","G","H","I","J") {
        $cell = $ws.Range($col + $r)
        if ($cell.Value2 -ne $null) {
            $cell.Font.Color = 255
        }
    }
}

# --- New trade row 7 (STT 6) added below the existing rows (no shifting of row 21 total) ---
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "17/12/2025"
$ws.Range("C7").Value = "mua"
$ws.Range("D7").Value = "btc/usdt "
$ws.Range("E7").Value = 88098
$ws.Range("F7").Value = 85000
$ws.Range("H7").Value = 261
$ws.Range("J7").Value = "bắt đáy , sai quy tắc vào lệnh "

foreach ($col in "A","B","C","D","E","F","H","J") {
    $ws.Range($col + "7").Font.Color = 255
}

# --- Selection moved to J11 ---
$ws.Range("J11").Select()
